$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2787.4644
$ws.Range("J17").Value = 2787.4644
$ws.Range("L17").Value = 8362.393199999999
$ws.Range("N17").Value = -8698.393199999999

$ws.Range("H41").Value = 3759.4
$ws.Range("J41").Value = 3699
$ws.Range("L41").Value = 3699
$ws.Range("N41").Value = -4579

$ws.Range("H51").Value = 5768.3335
$ws.Range("J51").Value = 5500
$ws.Range("L51").Value = 5500
$ws.Range("N51").Value = -6468

$ws.Range("H132").Value = 4543.7334
$ws.Range("I132").Value = 1396.6154
$ws.Range("K132").Value = 4189.8462
$ws.Range("M132").Value = -1659.8462

$ws.Range("H135").Value = 3082.75
$ws.Range("J135").Value = 3110.6667
$ws.Range("L135").Value = 27996.0003
$ws.Range("N135").Value = -33066.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9258.166999999999
$ws.Range("I74").Value = 8892.5
$ws.Range("K74").Value = 8892.5
$ws.Range("M74").Value = -8018.5

$ws.Range("H77").Value = 9258.166999999999
$ws.Range("I77").Value = 8892.5
$ws.Range("K77").Value = 44462.5
$ws.Range("M77").Value = -40094.5

$ws.Range("H88").Value = 1670
$ws.Range("J88").Value = 2107.5
$ws.Range("L88").Value = 2107.5
$ws.Range("N88").Value = -2919.5

$ws.Range("H91").Value = 1670
$ws.Range("J91").Value = 2107.5
$ws.Range("L91").Value = 2107.5
$ws.Range("N91").Value = -4915.5

$ws.Range("H101").Value = 37799.5
$ws.Range("J101").Value = 37799.5
$ws.Range("L101").Value = 37799.5
$ws.Range("N101").Value = -44289.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4670
$ws.Range("I86").Value = 2605.5
$ws.Range("J86").Value = 6734.5
$ws.Range("K86").Value = 2605.5
$ws.Range("L86").Value = 6734.5
$ws.Range("M86").Value = -1482.5
$ws.Range("N86").Value = -8980.5

$ws.Range("H89").Value = 4670
$ws.Range("I89").Value = 2605.5
$ws.Range("J89").Value = 6734.5
$ws.Range("K89").Value = 13027.5
$ws.Range("L89").Value = 33672.5
$ws.Range("M89").Value = -7411.5
$ws.Range("N89").Value = -44904.5

$ws.Range("H107").Value = 3145.5356
$ws.Range("I107").Value = 753.75
$ws.Range("J107").Value = 9125
$ws.Range("K107").Value = 753.75
$ws.Range("L107").Value = 9125
$ws.Range("M107").Value = 1166.25
$ws.Range("N107").Value = -12965

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H74").Value = 60805.668
$ws.Range("J74").Value = 60805.668
$ws.Range("L74").Value = 60805.668
$ws.Range("N74").Value = -62553.668

$ws.Range("H77").Value = 60805.668
$ws.Range("J77").Value = 60805.668
$ws.Range("L77").Value = 182417.004
$ws.Range("N77").Value = -191153.004

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H99").Value = 3002.125
$ws.Range("I99").Value = 2403.2
$ws.Range("J99").Value = 4000.3333
$ws.Range("K99").Value = 2403.2
$ws.Range("L99").Value = 4000.3333
$ws.Range("M99").Value = -905.1999999999998
$ws.Range("N99").Value = -6996.3333

$ws.Range("H116").Value = 80333.336
$ws.Range("J116").Value = 80333.336
$ws.Range("L116").Value = 80333.336
$ws.Range("N116").Value = -89511.336

$ws.Range("H126").Value = 3002.125
$ws.Range("I126").Value = 2403.2
$ws.Range("J126").Value = 4000.3333
$ws.Range("K126").Value = 7209.599999999999
$ws.Range("L126").Value = 12000.9999
$ws.Range("M126").Value = -4739.599999999999
$ws.Range("N126").Value = -16940.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 134.26666
$ws.Range("J23").Value = 133.81818
$ws.Range("L23").Value = 401.4545400000001
$ws.Range("N23").Value = -871.4545400000001

$ws.Range("H107").Value = 724.1429000000001
$ws.Range("J107").Value = 899.75
$ws.Range("L107").Value = 2699.25
$ws.Range("N107").Value = -6539.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4375.7085
$ws.Range("J46").Value = 4932.6665
$ws.Range("L46").Value = 4932.6665
$ws.Range("N46").Value = -5308.6665

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H93").Value = 1985.4445
$ws.Range("I93").Value = 1944.3334
$ws.Range("J93").Value = 2067.6667
$ws.Range("K93").Value = 1944.3334
$ws.Range("L93").Value = 2067.6667
$ws.Range("M93").Value = -696.3334
$ws.Range("N93").Value = -4563.6667

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H106").Value = 12898
$ws.Range("J106").Value = 12898
$ws.Range("L106").Value = 12898
$ws.Range("N106").Value = -15422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 111897.52
$ws.Range("I4").Value = 111897.52
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 111897.52
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -111784.52
$ws.Range("N4").ClearContents()

$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766

$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H117").Value = 70204.5
$ws.Range("J117").Value = 70204.5
$ws.Range("L117").Value = 70204.5
$ws.Range("N117").Value = -79382.5
